$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1: header / week-code labels (columns E..AE after the edit).
# The week numbering was renumbered (weeks 00,01,02,03,16 introduced; weeks
# 13,14,18,23,25,27 and everything past 34 dropped), shrinking the used
# range from A1:AL21 down to A1:AE21.
# ---------------------------------------------------------------------------
$headerValues = @{
    5  = "00-2020"   # E1
    6  = "01-2020"   # F1
    7  = "02-2020"   # G1
    8  = "03-2020"   # H1
    9  = "04-2020"   # I1
    10 = "06-2020"   # J1
    11 = "07-2020"   # K1
    12 = "09-2020"   # L1
    13 = "10-2020"   # M1
    14 = "11-2020"   # N1
    15 = "12-2020"   # O1
    16 = "15-2020"   # P1
    17 = "16-2020"   # Q1
    18 = "17-2020"   # R1
    19 = "19-2020"   # S1
    20 = "20-2020"   # T1
    21 = "21-2020"   # U1
    22 = "22-2020"   # V1
    23 = "24-2020"   # W1
    24 = "26-2020"   # X1
    25 = "28-2020"   # Y1
    26 = "29-2020"   # Z1
    27 = "30-2020"   # AA1
    28 = "31-2020"   # AB1
    29 = "32-2020"   # AC1
    30 = "33-2020"   # AD1
    31 = "34-2020"   # AE1
}
foreach ($col in $headerValues.Keys) {
    $ws.Cells.Item(1, $col).Value = $headerValues[$col]
}

# Columns AF..AL (32..38) no longer exist in the shrunk range - clear them.
for ($col = 32; $col -le 38; $col++) {
    $ws.Cells.Item(1, $col).Value = ""
}

# ---------------------------------------------------------------------------
# Row 2 (Art): newly-marked weeks.
# ---------------------------------------------------------------------------
$row2Add = @(5, 10, 11, 12, 13, 15, 16, 20, 21, 22, 25, 26, 29, 31)
foreach ($col in $row2Add) {
    $ws.Cells.Item(2, $col).Value = 1
}

# ---------------------------------------------------------------------------
# Row 6 (Dodgeball): old marks cleared, replaced by a single new one (AD6).
# ---------------------------------------------------------------------------
$row6Clear = @(5, 6, 8, 9, 10, 11, 14, 17, 20, 21, 22, 24, 29)
foreach ($col in $row6Clear) {
    $ws.Cells.Item(6, $col).Value = ""
}
$ws.Cells.Item(6, 30).Value = 1   # AD6

# ---------------------------------------------------------------------------
# Row 8 (Athletics): all marks removed.
# ---------------------------------------------------------------------------
$row8Clear = @(5, 6, 7, 9, 11, 12, 15, 16, 17, 20, 21, 22, 23, 24, 25, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36)
foreach ($col in $row8Clear) {
    $ws.Cells.Item(8, $col).Value = ""
}

# ---------------------------------------------------------------------------
# Row 15 (Choir): marks rearranged.
# ---------------------------------------------------------------------------
$row15Clear = @(7, 8, 9, 11, 14, 15, 16, 17, 23, 25, 26, 27, 28, 29, 33, 34, 35)
foreach ($col in $row15Clear) {
    $ws.Cells.Item(15, $col).Value = ""
}
$row15Add = @(5, 6, 7, 8, 9, 11, 12, 13, 14, 15, 17, 19, 20, 23, 24, 25, 27, 29, 30, 31)
foreach ($col in $row15Add) {
    $ws.Cells.Item(15, $col).Value = 1
}

# ---------------------------------------------------------------------------
# Row 17 (Gardening): one new mark.
# ---------------------------------------------------------------------------
$ws.Cells.Item(17, 31).Value = 1   # AE17

# ---------------------------------------------------------------------------
# Row 19 (Drama): all marks removed.
# ---------------------------------------------------------------------------
$row19Clear = @(5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38)
foreach ($col in $row19Clear) {
    $ws.Cells.Item(19, $col).Value = ""
}

# ---------------------------------------------------------------------------
# Row 21 (MFL): newly-marked weeks.
# ---------------------------------------------------------------------------
$row21Add = @(7, 10, 12, 15, 16, 18, 19, 21, 22, 27, 28, 31)
foreach ($col in $row21Add) {
    $ws.Cells.Item(21, $col).Value = 1
}

# ---------------------------------------------------------------------------
# Shrink the worksheet's used range back down to A1:AE21 by clearing
# anything that used to live in columns AF:AL across all 21 rows.
# ---------------------------------------------------------------------------
$ws.Range("AF1:AL21").Clear()
